$d = $word.ActiveDocument

$r = $d.Content
$r.Find.Execute(", Selection.activeObject")
$r.Collapse(0)
$r.Font.Size = 14
$r.InsertAfter(", OnGUI()")
$r.InsertAfter(", Handles.DrawBezier()..")
